$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = '2022-07-18 20:57:16'
$ws.Range("O3").Value = '2022-07-18 20:57:16'
$ws.Range("O4").Value = '2022-07-18 20:57:16'
$ws.Range("O5").Value = '2022-07-18 20:57:16'
$ws.Range("O6").Value = '2022-07-18 20:57:16'
$ws.Range("O7").Value = '2022-07-18 20:57:16'
$ws.Range("O8").Value = '2022-07-18 20:57:16'
$ws.Range("O9").Value = '2022-07-18 20:57:16'
$ws.Range("O10").Value = '2022-07-18 20:57:16'
$ws.Range("O11").Value = '2022-07-18 20:57:16'
$ws.Range("O12").Value = '2022-07-18 20:57:16'
$ws.Range("O13").Value = '2022-07-18 20:57:16'
$ws.Range("O14").Value = '2022-07-18 20:57:16'
$ws.Range("O15").Value = '2022-07-18 20:57:16'
$ws.Range("O16").Value = '2022-07-18 20:57:16'
$ws.Range("O17").Value = '2022-07-18 20:57:16'
$ws.Range("O18").Value = '2022-07-18 20:57:16'
$ws.Range("A19").Value = '''6761133'
$ws.Range("B19").Value = 'Duracell Batterien PLUS C/LR14 2 Stück'
$ws.Range("C19").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-clr14-2-stueck/p/6761133'
$ws.Range("D19").Value = '2ST'
$ws.Range("I19").Value = '4.98/1ST'
$ws.Range("K19").Value = '''4.98'
$ws.Range("M19").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Range("N19").Value = 'Duracell Batterien PLUS C/LR14 2 Stück 9.95 Schweizer Franken'
$ws.Range("O19").Value = '2022-07-18 20:57:16'
$ws.Range("A20").Value = '''6753557'
$ws.Range("B20").Value = 'Duracell Batterien PLUS AAA/LR03 4 Stück'
$ws.Range("C20").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-4-stueck/p/6753557'
$ws.Range("D20").Value = '4ST'
$ws.Range("I20").Value = '2.49/1ST'
$ws.Range("K20").Value = '''2.49'
$ws.Range("M20").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Range("N20").Value = 'Duracell Batterien PLUS AAA/LR03 4 Stück 9.95 Schweizer Franken'
$ws.Range("O20").Value = '2022-07-18 20:57:16'
$ws.Range("A21").Value = '''6761135'
$ws.Range("B21").Value = 'Duracell Batterie PLUS 9V/6LR61 1 Stück'
$ws.Range("C21").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-plus-9v6lr61-1-stueck/p/6761135'
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 'Duracell'
$ws.Range("H21").Value = '''9.95'
$ws.Range("I21").Value = '9.95/1ST'
$ws.Range("K21").Value = '''9.95'
$ws.Range("N21").Value = 'Duracell Batterie PLUS 9V/6LR61 1 Stück 9.95 Schweizer Franken'
$ws.Range("O21").Value = '2022-07-18 20:57:16'
$ws.Range("A22").Value = '''3494233'
$ws.Range("B22").Value = 'Varta Electronics CR2032 1er Bli'
$ws.Range("C22").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2032-1er-bli/p/3494233'
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 4.5
$ws.Range("G22").Value = 'Varta'
$ws.Range("H22").Value = '''4.95'
$ws.Range("I22").Value = '4.95/1ST'
$ws.Range("K22").Value = '''4.95'
$ws.Range("N22").Value = 'Varta Electronics CR2032 1er Bli 4.95 Schweizer Franken'
$ws.Range("O22").Value = '2022-07-18 20:57:16'
$ws.Range("O23").Value = '2022-07-18 20:57:16'
$ws.Range("O24").Value = '2022-07-18 20:57:16'
$ws.Range("O25").Value = '2022-07-18 20:57:16'
$ws.Range("O26").Value = '2022-07-18 20:57:16'
$ws.Range("O27").Value = '2022-07-18 20:57:16'
$ws.Range("A28").Value = '''6999781'
$ws.Range("B28").Value = 'Varta Batterien Longlife Power AAA/LR03 2x12 Stück'
$ws.Range("C28").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-batterien-longlife-power-aaalr03-2x12-stueck/p/6999781'
$ws.Range("D28").Value = '24ST'
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 'Varta'
$ws.Range("H28").Value = '''20.85'
$ws.Range("I28").Value = '0.87/1ST'
$ws.Range("K28").Value = '''0.87'
$ws.Range("M28").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Range("N28").Value = 'Varta Batterien Longlife Power AAA/LR03 2x12 Stück 50% Aktion 20.85 Schweizer Franken statt 41.70 Schweizer Franken'
$ws.Range("O28").Value = '2022-07-18 20:57:16'
$ws.Range("A29").Value = '''6999749'
$ws.Range("B29").Value = 'Varta Batterien Longlife Power AA/LR6 2x12 Stück'
$ws.Range("C29").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-batterien-longlife-power-aalr6-2x12-stueck/p/6999749'
$ws.Range("D29").Value = '24ST'
$ws.Range("G29").Value = 'Varta'
$ws.Range("H29").Value = '''20.85'
$ws.Range("I29").Value = '0.87/1ST'
$ws.Range("J29").Value = 'Preis pro 1 Stück'
$ws.Range("K29").Value = '''0.87'
$ws.Range("L29").Value = '1ST'
$ws.Range("M29").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Range("N29").Value = 'Varta Batterien Longlife Power AA/LR6 2x12 Stück 50% Aktion 20.85 Schweizer Franken statt 41.70 Schweizer Franken'
$ws.Range("O29").Value = '2022-07-18 20:57:16'
$ws.Range("A30").Value = '''6753554'
$ws.Range("B30").Value = 'Duracell Batterien Optimum AA/LR6 4 Stück'
$ws.Range("C30").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-optimum-aalr6-4-stueck/p/6753554'
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 'Duracell'
$ws.Range("H30").Value = '''11.95'
$ws.Range("I30").Value = '2.99/1ST'
$ws.Range("K30").Value = '''2.99'
$ws.Range("N30").Value = 'Duracell Batterien Optimum AA/LR6 4 Stück 11.95 Schweizer Franken'
$ws.Range("O30").Value = '2022-07-18 20:57:16'
$ws.Range("A31").Value = '''6508223'
$ws.Range("B31").Value = 'satrap Venti WS Standventilator'
$ws.Range("C31").Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-venti-ws-standventilator/p/6508223'
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 'satrap'
$ws.Range("H31").Value = '''59.95'
$ws.Range("I31").Value = ""
$ws.Range("J31").Value = ""
$ws.Range("K31").Value = ""
$ws.Range("L31").Value = ""
$ws.Range("M31").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Range("N31").Value = 'satrap Venti WS Standventilator 59.95 Schweizer Franken'
$ws.Range("O31").Value = '2022-07-18 20:57:16'
$ws.Range("A32").Value = '''3494130'
$ws.Range("B32").Value = 'Varta Longlife Power AA 4er Bli'
$ws.Range("C32").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-aa-4er-bli/p/3494130'
$ws.Range("D32").Value = '4ST'
$ws.Range("E32").Value = 1
$ws.Range("G32").Value = 'Varta'
$ws.Range("H32").Value = '''8.95'
$ws.Range("I32").Value = '2.24/1ST'
$ws.Range("K32").Value = '''2.24'
$ws.Range("M32").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Range("N32").Value = 'Varta Longlife Power AA 4er Bli 8.95 Schweizer Franken'
$ws.Range("O32").Value = '2022-07-18 20:57:16'
$ws.Range("A33").Value = '''3494065'
$ws.Range("B33").Value = 'Varta Longlife Max Power AA 4er Bli'
$ws.Range("C33").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-max-power-aa-4er-bli/p/3494065'
$ws.Range("D33").Value = '4ST'
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 5
$ws.Range("H33").Value = '''9.95'
$ws.Range("I33").Value = '2.49/1ST'
$ws.Range("K33").Value = '''2.49'
$ws.Range("M33").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Range("N33").Value = 'Varta Longlife Max Power AA 4er Bli 9.95 Schweizer Franken'
$ws.Range("O33").Value = '2022-07-18 20:57:16'
$ws.Range("A34").Value = '''5683906'
$ws.Range("B34").Value = 'Prix Garantie Wasserkocher 1l'
$ws.Range("C34").Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/prix-garantie-wasserkocher-1l/p/5683906'
$ws.Range("D34").Value = '1ST'
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 4
$ws.Range("H34").Value = '''19.95'
$ws.Range("I34").Value = '19.95/1ST'
$ws.Range("K34").Value = '''19.95'
$ws.Range("M34").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Range("N34").Value = 'Prix Garantie Wasserkocher 1l - Online kein Bestand 19.95 Schweizer Franken'
$ws.Range("O34").Value = '2022-07-18 20:57:16'
$ws.Range("A35").Value = '''4905484'
$ws.Range("B35").Value = 'Alkaline Batterie LR20/D 2 Stück'
$ws.Range("C35").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-lr20d-2-stueck/p/4905484'
$ws.Range("D35").Value = '2ST'
$ws.Range("E35").Value = 3
$ws.Range("F35").Value = 5
$ws.Range("G35").Value = 'Coop'
$ws.Range("H35").Value = '''5.95'
$ws.Range("I35").Value = '2.98/1ST'
$ws.Range("K35").Value = '''2.98'
$ws.Range("M35").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Range("N35").Value = 'Alkaline Batterie LR20/D 2 Stück 5.95 Schweizer Franken'
$ws.Range("O35").Value = '2022-07-18 20:57:16'
$ws.Range("O36").Value = '2022-07-18 20:57:16'
$ws.Range("O37").Value = '2022-07-18 20:57:16'
$ws.Range("O38").Value = '2022-07-18 20:57:16'
$ws.Range("O39").Value = '2022-07-18 20:57:16'
$ws.Range("O40").Value = '2022-07-18 20:57:16'
$ws.Range("O41").Value = '2022-07-18 20:57:16'
$ws.Range("A42").Value = '''6119284'
$ws.Range("B42").Value = 'satrap Aqua SA10 Wasserkocher'
$ws.Range("C42").Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-aqua-sa10-wasserkocher/p/6119284'
$ws.Range("E42").Value = 3
$ws.Range("F42").Value = 2.5
$ws.Range("H42").Value = '''49.95'
$ws.Range("N42").Value = 'satrap Aqua SA10 Wasserkocher 49.95 Schweizer Franken'
$ws.Range("O42").Value = '2022-07-18 20:57:16'
$ws.Range("A43").Value = '''5751576'
$ws.Range("B43").Value = 'satrap Toasty 1 Toaster'
$ws.Range("C43").Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-toasty-1-toaster/p/5751576'
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 5
$ws.Range("H43").Value = '''29.95'
$ws.Range("N43").Value = 'satrap Toasty 1 Toaster 29.95 Schweizer Franken'
$ws.Range("O43").Value = '2022-07-18 20:57:16'
$ws.Range("O44").Value = '2022-07-18 20:57:16'
$ws.Range("A45").Value = '''4358323'
$ws.Range("B45").Value = 'Rayovac Hörgerätebatterien 312 6 Stück'
$ws.Range("C45").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-312-6-stueck/p/4358323'
$ws.Range("D45").Value = '6ST'
$ws.Range("E45").Value = 1
$ws.Range("F45").Value = 4
$ws.Range("G45").Value = 'Rayovac'
$ws.Range("H45").Value = '''9.95'
$ws.Range("I45").Value = '1.66/1ST'
$ws.Range("K45").Value = '''1.66'
$ws.Range("N45").Value = 'Rayovac Hörgerätebatterien 312 6 Stück 9.95 Schweizer Franken'
$ws.Range("O45").Value = '2022-07-18 20:57:16'
$ws.Range("A46").Value = '''3494063'
$ws.Range("B46").Value = 'Varta Longlife Power C 2er Bli'
$ws.Range("C46").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-c-2er-bli/p/3494063'
$ws.Range("D46").Value = '2ST'
$ws.Range("E46").Value = ""
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 'Varta'
$ws.Range("H46").Value = '''7.95'
$ws.Range("I46").Value = '3.98/1ST'
$ws.Range("K46").Value = '''3.98'
$ws.Range("N46").Value = 'Varta Longlife Power C 2er Bli 7.95 Schweizer Franken'
$ws.Range("O46").Value = '2022-07-18 20:57:16'
$ws.Range("O47").Value = '2022-07-18 20:57:16'
$ws.Range("O48").Value = '2022-07-18 20:57:16'
$ws.Range("O49").Value = '2022-07-18 20:57:16'
$ws.Range("A50").Value = '''4942594'
$ws.Range("B50").Value = 'Skross Adapter World-Schweiz'
$ws.Range("C50").Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/skross-adapter-world-schweiz/p/4942594'
$ws.Range("D50").Value = ""
$ws.Range("E50").Value = ""
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 'Skross'
$ws.Range("H50").Value = '''19.95'
$ws.Range("I50").Value = ""
$ws.Range("J50").Value = ""
$ws.Range("K50").Value = ""
$ws.Range("L50").Value = ""
$ws.Range("M50").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Range("N50").Value = 'Skross Adapter World-Schweiz 19.95 Schweizer Franken'
$ws.Range("O50").Value = '2022-07-18 20:57:16'
$ws.Range("A51").Value = '''3494131'
$ws.Range("B51").Value = 'Varta Longlife Power AAA 4er Bli'
$ws.Range("C51").Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-aaa-4er-bli/p/3494131'
$ws.Range("D51").Value = '4ST'
$ws.Range("G51").Value = 'Varta'
$ws.Range("H51").Value = '''8.95'
$ws.Range("I51").Value = '2.24/1ST'
$ws.Range("J51").Value = 'Preis pro 1 Stück'
$ws.Range("K51").Value = '''2.24'
$ws.Range("L51").Value = '1ST'
$ws.Range("M51").Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Range("N51").Value = 'Varta Longlife Power AAA 4er Bli 8.95 Schweizer Franken'
$ws.Range("O51").Value = '2022-07-18 20:57:16'
$ws.Range("A52").Value = '''5763068'
$ws.Range("B52").Value = 'satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger'
$ws.Range("C52").Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-aspira-sine-a700-beutelloser-zyklonstaubsauger/p/5763068'
$ws.Range("E52").Value = 3
$ws.Range("F52").Value = 3.5
$ws.Range("G52").Value = 'satrap'
$ws.Range("H52").Value = '''119.00'
$ws.Range("N52").Value = 'satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger 119.00 Schweizer Franken'
$ws.Range("O52").Value = '2022-07-18 20:57:16'
$ws.Range("O53").Value = '2022-07-18 20:57:16'
$ws.Range("O54").Value = '2022-07-18 20:57:16'
$ws.Range("O55").Value = '2022-07-18 20:57:16'
$ws.Range("O56").Value = '2022-07-18 20:57:16'
$ws.Range("O57").Value = '2022-07-18 20:57:16'
$ws.Range("O58").Value = '2022-07-18 20:57:16'
$ws.Range("O59").Value = '2022-07-18 20:57:16'
$ws.Range("O60").Value = '2022-07-18 20:57:16'
$ws.Range("O61").Value = '2022-07-18 20:57:16'
$ws.Range("O62").Value = '2022-07-18 20:57:16'
$ws.Range("O63").Value = '2022-07-18 20:57:16'
$ws.Range("O64").Value = '2022-07-18 20:57:16'
$ws.Range("O65").Value = '2022-07-18 20:57:16'
$ws.Range("O66").Value = '2022-07-18 20:57:16'
$ws.Range("O67").Value = '2022-07-18 20:57:16'
$ws.Range("O68").Value = '2022-07-18 20:57:16'
$ws.Range("O69").Value = '2022-07-18 20:57:16'
$ws.Range("O70").Value = '2022-07-18 20:57:16'
$ws.Range("O71").Value = '2022-07-18 20:57:16'
$ws.Range("O72").Value = '2022-07-18 20:57:16'
$ws.Range("O73").Value = '2022-07-18 20:57:16'
$ws.Range("O74").Value = '2022-07-18 20:57:16'
$ws.Range("O75").Value = '2022-07-18 20:57:16'
$ws.Range("O76").Value = '2022-07-18 20:57:16'
$ws.Range("O77").Value = '2022-07-18 20:57:16'
$ws.Range("O78").Value = '2022-07-18 20:57:16'
$ws.Range("O79").Value = '2022-07-18 20:57:16'
$ws.Range("O80").Value = '2022-07-18 20:57:16'
$ws.Range("O81").Value = '2022-07-18 20:57:16'
$ws.Range("O82").Value = '2022-07-18 20:57:16'
$ws.Range("O83").Value = '2022-07-18 20:57:16'
$ws.Range("O84").Value = '2022-07-18 20:57:16'
$ws.Range("O85").Value = '2022-07-18 20:57:16'
$ws.Range("O86").Value = '2022-07-18 20:57:16'
$ws.Range("O87").Value = '2022-07-18 20:57:16'
$ws.Range("O88").Value = '2022-07-18 20:57:16'
